$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3945, 3985, 3985, 3985, 3985, 4047, 4047, 4125, 4385, 4385, 4474)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
